# Fruta / hortaliza, semanal
# Insert 3 new weekly rows at the top of the date-ordered block (rows 375-377),
# pushing the existing rows 375-392 down to 378-395, and fill the new rows
# with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 375 (shifts old rows 375:392 -> 378:395)
$ws.Range("A375:A377").EntireRow.Insert()

# New row 375
$ws.Cells.Item(375, 1).Value = 9
$ws.Cells.Item(375, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(375, 3).Value = "Metropolitana"
$ws.Cells.Item(375, 4).Value = 45267
$ws.Cells.Item(375, 5).Value = 13
$ws.Cells.Item(375, 6).Value = "Fruta"
$ws.Cells.Item(375, 7).Value = 100101
$ws.Cells.Item(375, 8).Value = "Berries"
$ws.Cells.Item(375, 9).Value = 100101001
$ws.Cells.Item(375, 10).Value = "Arándano (blue)"
$ws.Cells.Item(375, 11).Value = "Sin especificar"
$ws.Cells.Item(375, 12).Value = "Especial"
$ws.Cells.Item(375, 13).Value = 500
$ws.Cells.Item(375, 14).Value = 4000
$ws.Cells.Item(375, 15).Value = 4000
$ws.Cells.Item(375, 16).Value = 4000
$ws.Cells.Item(375, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(375, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(375, 19).Value = 2000
$ws.Cells.Item(375, 20).Value = 2

# New row 376
$ws.Cells.Item(376, 1).Value = 9
$ws.Cells.Item(376, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(376, 3).Value = "Metropolitana"
$ws.Cells.Item(376, 4).Value = 45267
$ws.Cells.Item(376, 5).Value = 13
$ws.Cells.Item(376, 6).Value = "Fruta"
$ws.Cells.Item(376, 7).Value = 100101
$ws.Cells.Item(376, 8).Value = "Berries"
$ws.Cells.Item(376, 9).Value = 100101001
$ws.Cells.Item(376, 10).Value = "Arándano (blue)"
$ws.Cells.Item(376, 11).Value = "Sin especificar"
$ws.Cells.Item(376, 12).Value = "Especial"
$ws.Cells.Item(376, 13).Value = 150
$ws.Cells.Item(376, 14).Value = 4000
$ws.Cells.Item(376, 15).Value = 4000
$ws.Cells.Item(376, 16).Value = 4000
$ws.Cells.Item(376, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(376, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(376, 19).Value = 2000
$ws.Cells.Item(376, 20).Value = 2

# New row 377
$ws.Cells.Item(377, 1).Value = 9
$ws.Cells.Item(377, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(377, 3).Value = "Metropolitana"
$ws.Cells.Item(377, 4).Value = 45267
$ws.Cells.Item(377, 5).Value = 13
$ws.Cells.Item(377, 6).Value = "Fruta"
$ws.Cells.Item(377, 7).Value = 100101
$ws.Cells.Item(377, 8).Value = "Berries"
$ws.Cells.Item(377, 9).Value = 100101001
$ws.Cells.Item(377, 10).Value = "Arándano (blue)"
$ws.Cells.Item(377, 11).Value = "Sin especificar"
$ws.Cells.Item(377, 12).Value = "Primera"
$ws.Cells.Item(377, 13).Value = 250
$ws.Cells.Item(377, 14).Value = 3600
$ws.Cells.Item(377, 15).Value = 3600
$ws.Cells.Item(377, 16).Value = 3600
$ws.Cells.Item(377, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(377, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(377, 19).Value = 1800
$ws.Cells.Item(377, 20).Value = 2
